$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised historical values (rows 388-390) ---
$ws.Range("B388").Value = 7431065591000
$ws.Range("D388").Value = 66946536855.85586

$ws.Range("B389").Value = 7779488501000
$ws.Range("D389").Value = 67465861599.16747

$ws.Range("B390").Value = 8162661859000
$ws.Range("D390").Value = 67914650628.17206

# --- Append new rows 399-401 ---
# Copy the formatting of the last existing data row (398) down into the
# three new rows so the new date cells (column A) keep the same style
# (bold, thin border, center/top alignment, date number format).
$ws.Range("A398:D398").Copy()
$ws.Range("A399:D401").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A399").Value = 44958
$ws.Range("B399").Value = 13551717200000
$ws.Range("C399").Value = 0.00507227998985544
$ws.Range("D399").Value = 68738103981.73979

$ws.Range("A400").Value = 44986
$ws.Range("B400").Value = 13436944380000
$ws.Range("C400").Value = 0.004785146904009953
$ws.Range("D400").Value = 64297752799.31094

$ws.Range("A401").Value = 45017
$ws.Range("B401").Value = 14573629300000
$ws.Range("C401").Value = 0.004491555874955084
$ws.Range("D401").Value = 65458270301.83256
